$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bus voltage magnitude results (380 kV case)
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.101054333566038
$ws.Cells.Item(2, 4).Value = 1.110752982992693
$ws.Cells.Item(2, 5).Value = 1.103164413499877
$ws.Cells.Item(2, 6).Value = 1.117972245795065
$ws.Cells.Item(2, 9).Value = 1.058028964873843
$ws.Cells.Item(2, 10).Value = 1.10583378829899
$ws.Cells.Item(2, 11).Value = 1.11334950570567
$ws.Cells.Item(2, 12).Value = 1.105779761942219
$ws.Cells.Item(2, 13).Value = 1.120551114044796
$ws.Cells.Item(2, 14).Value = 1.107404199119283
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.102989792001522
$ws.Cells.Item(3, 4).Value = 1.112635034171671
$ws.Cells.Item(3, 5).Value = 1.104939490979827
$ws.Cells.Item(3, 6).Value = 1.119907604759193
$ws.Cells.Item(3, 9).Value = 1.058533279833621
$ws.Cells.Item(3, 10).Value = 1.107433775406764
$ws.Cells.Item(3, 11).Value = 1.115052904158254
$ws.Cells.Item(3, 12).Value = 1.107375049459709
$ws.Cells.Item(3, 13).Value = 1.122308997660593
$ws.Cells.Item(3, 14).Value = 1.109006458392271
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.104238409224084
$ws.Cells.Item(4, 4).Value = 1.113849365977964
$ws.Cells.Item(4, 5).Value = 1.106084435314537
$ws.Cells.Item(4, 6).Value = 1.121156466992827
$ws.Cells.Item(4, 9).Value = 1.058856345452827
$ws.Cells.Item(4, 10).Value = 1.108464993262522
$ws.Cells.Item(4, 11).Value = 1.116151187825904
$ws.Cells.Item(4, 12).Value = 1.108403193068623
$ws.Cells.Item(4, 13).Value = 1.123442580197824
$ws.Cells.Item(4, 14).Value = 1.110039140695666
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.104762449795251
$ws.Cells.Item(5, 4).Value = 1.114359056991372
$ws.Cells.Item(5, 5).Value = 1.106564914903571
$ws.Cells.Item(5, 6).Value = 1.121680683974865
$ws.Cells.Item(5, 9).Value = 1.058991388265073
$ws.Cells.Item(5, 10).Value = 1.108897557669275
$ws.Cells.Item(5, 11).Value = 1.116611983837083
$ws.Cells.Item(5, 12).Value = 1.108834457342321
$ws.Cells.Item(5, 13).Value = 1.123918226783871
$ws.Cells.Item(5, 14).Value = 1.110472319393492
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.104850387613674
$ws.Cells.Item(6, 4).Value = 1.114444589174882
$ws.Cells.Item(6, 5).Value = 1.106645539942
$ws.Cells.Item(6, 6).Value = 1.121768655699227
$ws.Cells.Item(6, 9).Value = 1.059014017347132
$ws.Cells.Item(6, 10).Value = 1.108970131367547
$ws.Cells.Item(6, 11).Value = 1.116689299947007
$ws.Cells.Item(6, 12).Value = 1.108906812297229
$ws.Cells.Item(6, 13).Value = 1.123998037027761
$ws.Cells.Item(6, 14).Value = 1.110544996154741
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.10424541490229
$ws.Cells.Item(7, 4).Value = 1.113856179666739
$ws.Cells.Item(7, 5).Value = 1.106090858840524
$ws.Cells.Item(7, 6).Value = 1.121163474740673
$ws.Cells.Item(7, 9).Value = 1.058858152933874
$ws.Cells.Item(7, 10).Value = 1.108470776956171
$ws.Cells.Item(7, 11).Value = 1.116157348604616
$ws.Cells.Item(7, 12).Value = 1.10840895941983
$ws.Cells.Item(7, 13).Value = 1.123448939366943
$ws.Cells.Item(7, 14).Value = 1.110044932602824
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.101709220492709
$ws.Cells.Item(8, 4).Value = 1.11138976348219
$ws.Cells.Item(8, 5).Value = 1.103765076198789
$ws.Cells.Item(8, 6).Value = 1.11862703464875
$ws.Cells.Item(8, 9).Value = 1.058200079359981
$ws.Cells.Item(8, 10).Value = 1.106375366896248
$ws.Cells.Item(8, 11).Value = 1.113926001715994
$ws.Cells.Item(8, 12).Value = 1.106319759336098
$ws.Cells.Item(8, 13).Value = 1.121146014936928
$ws.Cells.Item(8, 14).Value = 1.107946546820269
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.097210444145246
$ws.Cells.Item(9, 4).Value = 1.107016089428074
$ws.Cells.Item(9, 5).Value = 1.099637955867769
$ws.Cells.Item(9, 6).Value = 1.114130224553055
$ws.Cells.Item(9, 9).Value = 1.057015211030872
$ws.Cells.Item(9, 10).Value = 1.102650948272683
$ws.Cells.Item(9, 11).Value = 1.109963171007973
$ws.Cells.Item(9, 12).Value = 1.102606017408446
$ws.Cells.Item(9, 13).Value = 1.117057362242372
$ws.Cells.Item(9, 14).Value = 1.104216839095061
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.094190001397307
$ws.Cells.Item(10, 4).Value = 1.104080568257523
$ws.Cells.Item(10, 5).Value = 1.096865996989373
$ws.Cells.Item(10, 6).Value = 1.111112747520449
$ws.Cells.Item(10, 9).Value = 1.056207904682007
$ws.Cells.Item(10, 10).Value = 1.100145358318824
$ws.Cells.Item(10, 11).Value = 1.107299345823294
$ws.Cells.Item(10, 12).Value = 1.100107351346437
$ws.Cells.Item(10, 13).Value = 1.114309827048686
$ws.Cells.Item(10, 14).Value = 1.101707690916073
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.092876781716724
$ws.Cells.Item(11, 4).Value = 1.102804498595107
$ws.Cells.Item(11, 5).Value = 1.095660569430133
$ws.Cells.Item(11, 6).Value = 1.109801213914177
$ws.Cells.Item(11, 9).Value = 1.055854113793325
$ws.Cells.Item(11, 10).Value = 1.09905479187921
$ws.Cells.Item(11, 11).Value = 1.106140417918077
$ws.Cells.Item(11, 12).Value = 1.099019733826887
$ws.Cells.Item(11, 13).Value = 1.113114684660837
$ws.Cells.Item(11, 14).Value = 1.100615575747025
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.09238816464194
$ws.Cells.Item(12, 4).Value = 1.102329739052891
$ws.Cells.Item(12, 5).Value = 1.095212023050006
$ws.Cells.Item(12, 6).Value = 1.109313284315604
$ws.Cells.Item(12, 9).Value = 1.055722057722872
$ws.Cells.Item(12, 10).Value = 1.098648839722484
$ws.Cells.Item(12, 11).Value = 1.105709095644322
$ws.Cells.Item(12, 12).Value = 1.098614869376355
$ws.Cells.Item(12, 13).Value = 1.112669914933994
$ws.Cells.Item(12, 14).Value = 1.100209047091673
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.092493012579027
$ws.Cells.Item(13, 4).Value = 1.102431611836777
$ws.Cells.Item(13, 5).Value = 1.095308274206025
$ws.Cells.Item(13, 6).Value = 1.109417982004422
$ws.Cells.Item(13, 9).Value = 1.055750413398475
$ws.Cells.Item(13, 10).Value = 1.098735957429474
$ws.Cells.Item(13, 11).Value = 1.10580165433083
$ws.Cells.Item(13, 12).Value = 1.098701754117295
$ws.Cells.Item(13, 13).Value = 1.112765357964704
$ws.Cells.Item(13, 14).Value = 1.1002962885158
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.092836409532227
$ws.Cells.Item(14, 4).Value = 1.102765270679212
$ws.Cells.Item(14, 5).Value = 1.095623508837494
$ws.Cells.Item(14, 6).Value = 1.109760897303611
$ws.Cells.Item(14, 9).Value = 1.055843211148952
$ws.Cells.Item(14, 10).Value = 1.099021253563599
$ws.Cells.Item(14, 11).Value = 1.106104782049928
$ws.Cells.Item(14, 12).Value = 1.098986285579083
$ws.Cells.Item(14, 13).Value = 1.113077937139841
$ws.Cells.Item(14, 14).Value = 1.100581989803159
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.093047877138991
$ws.Cells.Item(15, 4).Value = 1.102970746093517
$ws.Cells.Item(15, 5).Value = 1.095817629005289
$ws.Cells.Item(15, 6).Value = 1.109972076280647
$ws.Cells.Item(15, 9).Value = 1.055900301539929
$ws.Cells.Item(15, 10).Value = 1.099196918365189
$ws.Cells.Item(15, 11).Value = 1.10629143641502
$ws.Cells.Item(15, 12).Value = 1.099161478218641
$ws.Cells.Item(15, 13).Value = 1.113270415397882
$ws.Cells.Item(15, 14).Value = 1.100757904068916
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.094277040631013
$ws.Cells.Item(16, 4).Value = 1.104165150068009
$ws.Cells.Item(16, 5).Value = 1.096945886761176
$ws.Cells.Item(16, 6).Value = 1.111199683391564
$ws.Cells.Item(16, 9).Value = 1.05623129497966
$ws.Cells.Item(16, 10).Value = 1.100217615190338
$ws.Cells.Item(16, 11).Value = 1.107376142811817
$ws.Cells.Item(16, 12).Value = 1.100179411449052
$ws.Cells.Item(16, 13).Value = 1.114389028101345
$ws.Cells.Item(16, 14).Value = 1.101780050400633
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.095046612816651
$ws.Cells.Item(17, 4).Value = 1.104913020898275
$ws.Cells.Item(17, 5).Value = 1.097652217797389
$ws.Cells.Item(17, 6).Value = 1.111968387469877
$ws.Cells.Item(17, 9).Value = 1.056437782421467
$ws.Cells.Item(17, 10).Value = 1.100856349117276
$ws.Cells.Item(17, 11).Value = 1.108055069179867
$ws.Cells.Item(17, 12).Value = 1.100816398492635
$ws.Cells.Item(17, 13).Value = 1.115089231234418
$ws.Cells.Item(17, 14).Value = 1.102419691403011
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.095494976972012
$ws.Cells.Item(18, 4).Value = 1.105348763428927
$ws.Cells.Item(18, 5).Value = 1.09806371330695
$ws.Cells.Item(18, 6).Value = 1.112416283994993
$ws.Cells.Item(18, 9).Value = 1.05655781625182
$ws.Cells.Item(18, 10).Value = 1.101228370113633
$ws.Cells.Item(18, 11).Value = 1.108450548810761
$ws.Cells.Item(18, 12).Value = 1.101187395835013
$ws.Cells.Item(18, 13).Value = 1.115497124423706
$ws.Cells.Item(18, 14).Value = 1.102792240711854
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.095647771092851
$ws.Cells.Item(19, 4).Value = 1.105497260032057
$ws.Cells.Item(19, 5).Value = 1.098203939299688
$ws.Cells.Item(19, 6).Value = 1.112568925209536
$ws.Cells.Item(19, 9).Value = 1.056598675928403
$ws.Cells.Item(19, 10).Value = 1.10135512842342
$ws.Cells.Item(19, 11).Value = 1.108585308493152
$ws.Cells.Item(19, 12).Value = 1.101313804311687
$ws.Cells.Item(19, 13).Value = 1.115636117368365
$ws.Cells.Item(19, 14).Value = 1.10291917903298
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.094964098381255
$ws.Cells.Item(20, 4).Value = 1.104832831036144
$ws.Cells.Item(20, 5).Value = 1.097576486572795
$ws.Cells.Item(20, 6).Value = 1.111885962163141
$ws.Cells.Item(20, 9).Value = 1.056415670408737
$ws.Cells.Item(20, 10).Value = 1.10078787519671
$ws.Cells.Item(20, 11).Value = 1.107982281439688
$ws.Cells.Item(20, 12).Value = 1.100748112486936
$ws.Cells.Item(20, 13).Value = 1.115014160374452
$ws.Cells.Item(20, 14).Value = 1.102351120241624
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.092735310772911
$ws.Cells.Item(21, 4).Value = 1.102667037927515
$ws.Cells.Item(21, 5).Value = 1.09553070228065
$ws.Cells.Item(21, 6).Value = 1.109659938683581
$ws.Cells.Item(21, 9).Value = 1.055815902322775
$ws.Cells.Item(21, 10).Value = 1.098937265076333
$ws.Cells.Item(21, 11).Value = 1.106015541995572
$ws.Cells.Item(21, 12).Value = 1.09890252248218
$ws.Cells.Item(21, 13).Value = 1.112985913756183
$ws.Cells.Item(21, 14).Value = 1.100497882042607
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.091329175483493
$ws.Cells.Item(22, 4).Value = 1.101300848205046
$ws.Cells.Item(22, 5).Value = 1.094239814179975
$ws.Cells.Item(22, 6).Value = 1.108255895767877
$ws.Cells.Item(22, 9).Value = 1.055435083607226
$ws.Cells.Item(22, 10).Value = 1.097768684712124
$ws.Cells.Item(22, 11).Value = 1.104774075479113
$ws.Cells.Item(22, 12).Value = 1.097737054049719
$ws.Cells.Item(22, 13).Value = 1.111705799309234
$ws.Cells.Item(22, 14).Value = 1.099327642160243
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.092075058699898
$ws.Cells.Item(23, 4).Value = 1.102025522932857
$ws.Cells.Item(23, 5).Value = 1.094924584294969
$ws.Cells.Item(23, 6).Value = 1.109000635904519
$ws.Cells.Item(23, 9).Value = 1.05563731822517
$ws.Cells.Item(23, 10).Value = 1.09838865514216
$ws.Cells.Item(23, 11).Value = 1.105432672342628
$ws.Cells.Item(23, 12).Value = 1.09835537908425
$ws.Cells.Item(23, 13).Value = 1.112384882114025
$ws.Cells.Item(23, 14).Value = 1.099948493019402
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.095001384691263
$ws.Cells.Item(24, 4).Value = 1.104869066859125
$ws.Cells.Item(24, 5).Value = 1.097610707783665
$ws.Cells.Item(24, 6).Value = 1.111923208080143
$ws.Cells.Item(24, 9).Value = 1.056425663132297
$ws.Cells.Item(24, 10).Value = 1.100818817287264
$ws.Cells.Item(24, 11).Value = 1.108015172713072
$ws.Cells.Item(24, 12).Value = 1.100778969681433
$ws.Cells.Item(24, 13).Value = 1.11504808328154
$ws.Cells.Item(24, 14).Value = 1.102382106273495
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.098377135963571
$ws.Cells.Item(25, 4).Value = 1.108150177193572
$ws.Cells.Item(25, 5).Value = 1.100708448826841
$ws.Cells.Item(25, 6).Value = 1.115296119467919
$ws.Cells.Item(25, 9).Value = 1.057324563192102
$ws.Cells.Item(25, 10).Value = 1.103617709203385
$ws.Cells.Item(25, 11).Value = 1.110991438948261
$ws.Cells.Item(25, 12).Value = 1.103570051244677
$ws.Cells.Item(25, 13).Value = 1.11811812540738
$ws.Cells.Item(25, 14).Value = 1.105184972937174
